# Commit: "Sua lai ket qua hoc tap va sinh vien"
# Insert a new worksheet "KetQuaHocTap" right after "HocPhan" (before "MonHoc"),
# populate it with the exam-result data, and tweak a couple of view settings
# (active sheet / first visible tab, and the LopHocPhan sheet scroll position)
# to match the target workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new sheet right after "HocPhan" (so it lands before "MonHoc")
# ---------------------------------------------------------------------------
$hocPhan = $wb.Worksheets.Item("HocPhan")
$ws = $wb.Worksheets.Add($null, $hocPhan)
$ws.Name = "KetQuaHocTap"

# ---------------------------------------------------------------------------
# 2) Headers
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Điểm CK"
$ws.Range("B1").Value = "Điểm GK"
$ws.Range("C1").Value = "TH1"
$ws.Range("D1").Value = "TH2"
$ws.Range("E1").Value = "TK1"
$ws.Range("F1").Value = "TK2"
$ws.Range("G1").Value = "TK3"
$ws.Range("H1").Value = "mã lớp hp"
$ws.Range("I1").Value = "mã sinh viên"

# ---------------------------------------------------------------------------
# 3) Data rows (2..10)
# ---------------------------------------------------------------------------
$data = @(
    @(8,  8.5, 8,  8,  8,  8,  8,  420000012, "18000003"),
    @(4,  4,   4,  4,  4,  4,  4,  420000013, "18000002"),
    @(10, 10,  10, 10, 10, 10, 10, 420000013, "18000004"),
    @(8,  8.5, 8,  8,  8,  8,  8,  420000012, "18000001"),
    @(4,  4,   4,  4,  4,  4,  4,  420000013, "18000001"),
    @(10, 10,  10, 10, 10, 10, 10, 420000011, "18000001"),
    @(8,  8.5, 8,  8,  8,  8,  8,  420000008, "18000001"),
    @(4,  4,   4,  4,  4,  4,  4,  420000009, "18000001"),
    @(10, 10,  10, 10, 10, 10, 10, 420000010, "18000001")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = "'" + $r[8]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 4) Formatting
#    Order matters: the quote-prefixed / wrap-text style for column H must be
#    created before the 2-decimal numeric style for columns A:G, so that the
#    new style table entries land at indexes 11 and 12 respectively.
# ---------------------------------------------------------------------------
$colH = $ws.Range("H2:H10")
$colH.WrapText = $true
$colH.VerticalAlignment = -4108

$colI = $ws.Range("I2:I10")
$colI.VerticalAlignment = $null

$gradeBlock = $ws.Range("A1:G10")
$gradeBlock.NumberFormat = "0.00"

# Cells that were left on the General format (no decimals) instead of 0.00
$ws.Range("A3").NumberFormat = "General"
$ws.Range("A4:G4").NumberFormat = "General"
$ws.Range("A6").NumberFormat = "General"
$ws.Range("A7:G7").NumberFormat = "General"
$ws.Range("A9").NumberFormat = "General"
$ws.Range("A10:G10").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 5) Column widths / sheet layout
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19
$ws.Columns.Item(2).ColumnWidth = 7.9
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(4).ColumnWidth = 4.7
$ws.Columns.Item(5).ColumnWidth = 4.7
$ws.Columns.Item(6).ColumnWidth = 4.6
$ws.Columns.Item(7).ColumnWidth = 4.6
$ws.Columns.Item(8).ColumnWidth = 11
$ws.Columns.Item(9).ColumnWidth = 14.2

$ws.Range("M10").Select()

# ---------------------------------------------------------------------------
# 6) LopHocPhan: drop tabSelected, move scroll position from D1 to D7
# ---------------------------------------------------------------------------
$lopHocPhan = $wb.Worksheets.Item("LopHocPhan")
$lopHocPhan.Select()
$excel.ActiveWindow.ScrollColumn = 7
$lopHocPhan.Range("H12").Select()

# ---------------------------------------------------------------------------
# 7) Make KetQuaHocTap the active sheet / tab, with HocPhan as first visible
# ---------------------------------------------------------------------------
$ws.Select()
$excel.ActiveWindow.ScrollWorkbookTabs(0)
$hocPhan2 = $wb.Worksheets.Item("HocPhan")
$excel.ActiveWindow.ScrollWorkbookTabs($hocPhan2.Index - 1)
$ws.Select()
$ws.Range("M10").Select()

# ---------------------------------------------------------------------------
# 8) Window geometry
# ---------------------------------------------------------------------------
$excel.ActiveWindow.WindowState = -4143
$excel.Left = -108
$excel.Top = -108
$excel.Width = 23256
$excel.Height = 12576
